$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Knee Forward
$ws.Range("C2").Value = -0.3212
$ws.Range("E2").Value = 42.2784

# Row 3 - Knee Backward
$ws.Range("B3").Value = 0.71
$ws.Range("C3").Value = 1.0570999999999999
$ws.Range("D3").Value = 109.5784
$ws.Range("E3").Value = 5.0027999999999997

# Update the selection to match the diff (selection moved to G16)
$ws.Range("G16").Select()
